# Auto-generated edit script: refreshes cached market-price / profit figures
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets (scheduled market-data runner).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 57694068
$ws.Range("I135").Value = 20835114
$ws.Range("J135").Value = 500001500
$ws.Range("K135").Value = 187516026
$ws.Range("L135").Value = 4500013500
$ws.Range("M135").Value = -187513491
$ws.Range("N135").Value = -4500018570
$ws.Range("H137").Value = 3886
$ws.Range("I137").Value = 2078
$ws.Range("J137").Value = 5920
$ws.Range("K137").Value = 6234
$ws.Range("L137").Value = 17760
$ws.Range("M137").Value = -3684
$ws.Range("N137").Value = -22860
$ws.Range("H138").Value = 2394.2727
$ws.Range("I138").Value = 992.5
$ws.Range("J138").Value = 3276.8704
$ws.Range("K138").Value = 2977.5
$ws.Range("L138").Value = 9830.611199999999
$ws.Range("M138").Value = 2162.5
$ws.Range("N138").Value = -20110.6112

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 8146.619
$ws.Range("I61").Value = 3925.1724
$ws.Range("J61").Value = 17563.691
$ws.Range("K61").Value = 3925.1724
$ws.Range("L61").Value = 17563.691
$ws.Range("M61").Value = -3713.1724
$ws.Range("N61").Value = -17987.691
$ws.Range("H132").Value = 5150.6978
$ws.Range("I132").Value = 1961.95
$ws.Range("J132").Value = 7923.522
$ws.Range("K132").Value = 5885.85
$ws.Range("L132").Value = 23770.566
$ws.Range("M132").Value = -3355.85
$ws.Range("N132").Value = -28830.566
$ws.Range("H136").Value = 8146.619
$ws.Range("I136").Value = 3925.1724
$ws.Range("J136").Value = 17563.691
$ws.Range("K136").Value = 11775.5172
$ws.Range("L136").Value = 52691.073
$ws.Range("M136").Value = -9225.5172
$ws.Range("N136").Value = -57791.073

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H111").Value = 59100
$ws.Range("J111").Value = 59100
$ws.Range("L111").Value = 59100
$ws.Range("N111").Value = -67280
$ws.Range("H114").Value = 42684
$ws.Range("J114").Value = 42684
$ws.Range("L114").Value = 42684
$ws.Range("N114").Value = -51362
$ws.Range("H134").Value = 44381.918
$ws.Range("I134").Value = 2877.1667
$ws.Range("J134").Value = 168896.17
$ws.Range("K134").Value = 8631.500100000001
$ws.Range("L134").Value = 506688.51
$ws.Range("M134").Value = -6096.500100000001
$ws.Range("N134").Value = -511758.51

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3405.7646
$ws.Range("I31").Value = 2558.2727
$ws.Range("J31").Value = 4959.5
$ws.Range("K31").Value = 2558.2727
$ws.Range("L31").Value = 4959.5
$ws.Range("M31").Value = -2263.2727
$ws.Range("N31").Value = -5549.5
$ws.Range("H34").Value = 3405.7646
$ws.Range("I34").Value = 2558.2727
$ws.Range("J34").Value = 4959.5
$ws.Range("K34").Value = 2558.2727
$ws.Range("L34").Value = 4959.5
$ws.Range("M34").Value = -2356.2727
$ws.Range("N34").Value = -5363.5
$ws.Range("H41").Value = 23216.5
$ws.Range("I41").Value = 9999.5
$ws.Range("J41").Value = 29825
$ws.Range("K41").Value = 9999.5
$ws.Range("L41").Value = 29825
$ws.Range("M41").Value = -9571.5
$ws.Range("N41").Value = -30681
$ws.Range("H50").Value = 21582.5
$ws.Range("J50").Value = 21582.5
$ws.Range("L50").Value = 21582.5
$ws.Range("N50").Value = -22832.5
$ws.Range("H51").Value = 27191.615
$ws.Range("J51").Value = 27191.615
$ws.Range("L51").Value = 27191.615
$ws.Range("N51").Value = -28663.615
$ws.Range("H58").Value = 4333466.5
$ws.Range("I58").Value = 8268454
$ws.Range("J58").Value = 4980
$ws.Range("K58").Value = 8268454
$ws.Range("L58").Value = 4980
$ws.Range("M58").Value = -8268251
$ws.Range("N58").Value = -5386
$ws.Range("H59").Value = 17999
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H60").Value = 12933
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()
$ws.Range("H61").Value = 27191.615
$ws.Range("J61").Value = 27191.615
$ws.Range("L61").Value = 27191.615
$ws.Range("N61").Value = -27887.615
$ws.Range("H68").Value = 35492
$ws.Range("J68").Value = 35492
$ws.Range("L68").Value = 35492
$ws.Range("N68").Value = -36990
$ws.Range("H71").Value = 35492
$ws.Range("J71").Value = 35492
$ws.Range("L71").Value = 106476
$ws.Range("N71").Value = -113964
$ws.Range("H111").Value = 78000.336
$ws.Range("J111").Value = 78000.336
$ws.Range("L111").Value = 78000.336
$ws.Range("N111").Value = -86180.336
$ws.Range("H136").Value = 4333466.5
$ws.Range("I136").Value = 8268454
$ws.Range("J136").Value = 4980
$ws.Range("K136").Value = 24805362
$ws.Range("L136").Value = 14940
$ws.Range("M136").Value = -24802812
$ws.Range("N136").Value = -20040
$ws.Range("H141").Value = 39813.047
$ws.Range("J141").Value = 40409.367
$ws.Range("L141").Value = 40409.367
$ws.Range("N141").Value = -50769.367

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 30346212
$ws.Range("I9").Value = 999
$ws.Range("J9").Value = 33380734
$ws.Range("K9").Value = 2997
$ws.Range("L9").Value = 100142202
$ws.Range("M9").Value = -2773
$ws.Range("N9").Value = -100142650
$ws.Range("H108").Value = 1696.1428
$ws.Range("I108").Value = 808.1667
$ws.Range("J108").Value = 7024
$ws.Range("K108").Value = 2424.5001
$ws.Range("L108").Value = 21072
$ws.Range("M108").Value = 455.4998999999998
$ws.Range("N108").Value = -26832
$ws.Range("H114").Value = 361.30768
$ws.Range("I114").Value = 308.08334
$ws.Range("J114").Value = 1000
$ws.Range("K114").Value = 924.2500200000001
$ws.Range("L114").Value = 3000
$ws.Range("M114").Value = 2329.74998
$ws.Range("N114").Value = -9508
$ws.Range("H129").Value = 1497.5927
$ws.Range("I129").Value = 1056.4375
$ws.Range("J129").Value = 2139.2727
$ws.Range("K129").Value = 3169.3125
$ws.Range("L129").Value = 6417.8181
$ws.Range("M129").Value = 1830.6875
$ws.Range("N129").Value = -16417.8181
$ws.Range("H132").Value = 1724.9333
$ws.Range("I132").Value = 1499.8572
$ws.Range("J132").Value = 1921.875
$ws.Range("K132").Value = 13498.7148
$ws.Range("L132").Value = 17296.875
$ws.Range("M132").Value = -10968.7148
$ws.Range("N132").Value = -22356.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 14963.5
$ws.Range("J48").Value = 15000
$ws.Range("L48").Value = 15000
$ws.Range("N48").Value = -15970
$ws.Range("H59").Value = 34305
$ws.Range("J59").Value = 34305
$ws.Range("L59").Value = 34305
$ws.Range("N59").Value = -35471
$ws.Range("H102").Value = 2693.9805
$ws.Range("I102").Value = 2485.9167
$ws.Range("K102").Value = 2485.9167
$ws.Range("M102").Value = -863.9167000000002
$ws.Range("H132").Value = 9440.679
$ws.Range("I132").Value = 15305.875
$ws.Range("J132").Value = 7094.6
$ws.Range("K132").Value = 45917.625
$ws.Range("L132").Value = 21283.8
$ws.Range("M132").Value = -43387.625
$ws.Range("N132").Value = -26343.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1134.5
$ws.Range("I16").Value = 991
$ws.Range("K16").Value = 991
$ws.Range("M16").Value = -821
$ws.Range("H46").Value = 808.38464
$ws.Range("J46").Value = 849.9
$ws.Range("L46").Value = 849.9
$ws.Range("N46").Value = -1225.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H114").Value = 39800
$ws.Range("J114").Value = 39800
$ws.Range("L114").Value = 39800
$ws.Range("N114").Value = -48478
$ws.Range("H136").Value = 5619.5454
$ws.Range("I136").Value = 2533.4783
$ws.Range("J136").Value = 8999.522999999999
$ws.Range("K136").Value = 7600.4349
$ws.Range("L136").Value = 26998.569
$ws.Range("M136").Value = -5050.4349
$ws.Range("N136").Value = -32098.569
